# captaintsubasapy - using opencv for template matching
# Adds a "Description" column to Sheet1, populates four new template rows,
# and adds a new "Sheet2" that documents the "Location X" / "Location Y"
# columns used by the template-matching helper script.

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Sheet1
# ---------------------------------------------------------------------

# Make the header row bold (A1 keeps its text-number format + center
# alignment, B1 keeps the plain center alignment) - this creates the two
# new bold styles used throughout the sheet.
$ws1.Range("A1").Font.Bold = $true
$ws1.Range("B1").Font.Bold = $true

# New header cell, styled like the rest of the (now bold) header row.
$ws1.Range("F1").Value = "Description"
$ws1.Range("B1").Copy()
$ws1.Range("C1:F1").PasteSpecial($xlPasteFormats)

# Refresh the sample coordinates for the existing "001" row and give it
# a description.
$ws1.Range("B2").Value = 314
$ws1.Range("C2").Value = 141
$ws1.Range("D2").Value = 467
$ws1.Range("E2").Value = 303
$ws1.Range("F2").Value = "Run App"
$ws1.Range("B2").Copy()
$ws1.Range("F2").PasteSpecial($xlPasteFormats)

# New row: 002 - Enter App
$ws1.Range("A3").Value = "002"
$ws1.Range("B3").Value = 690
$ws1.Range("C3").Value = 479
$ws1.Range("D3").Value = 891
$ws1.Range("E3").Value = 631
$ws1.Range("F3").Value = "Enter App"

# New row: 003 - Go to Story Mode
$ws1.Range("A4").Value = "003"
$ws1.Range("B4").Value = 1029
$ws1.Range("C4").Value = 465
$ws1.Range("D4").Value = 1423
$ws1.Range("E4").Value = 565
$ws1.Range("F4").Value = "Go to Story Mode"

# New row: 004 - Go to Story Mode - Second Page
$ws1.Range("A5").Value = "004"
$ws1.Range("B5").Value = 203
$ws1.Range("C5").Value = 726
$ws1.Range("D5").Value = 499
$ws1.Range("E5").Value = 822
$ws1.Range("F5").Value = "Go to Story Mode - Second Page"

# New row: 005 - Middle School Part 1/3
$ws1.Range("A6").Value = "005"
$ws1.Range("B6").Value = 227
$ws1.Range("C6").Value = 306
$ws1.Range("D6").Value = 525
$ws1.Range("E6").Value = 502
$ws1.Range("F6").Value = "Middle School Part 1/3"

# Copy the already-established column styles down into the new rows.
$ws1.Range("A2").Copy()
$ws1.Range("A3:A6").PasteSpecial($xlPasteFormats)
$ws1.Range("B2").Copy()
$ws1.Range("B3:F6").PasteSpecial($xlPasteFormats)

# Widen the new description column.
$ws1.Columns.Item(6).ColumnWidth = 61.25

# Sheet1 is no longer the active/selected sheet once Sheet2 is added;
# its last selection moves to F1.
$null = $ws1.Range("F1").Select()

# ---------------------------------------------------------------------
# Sheet2: new sheet describing the Location X / Location Y columns.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Location X"
$ws2.Range("B1").Value = "Location Y"
$ws2.Range("C1").Value = "Description"

$ws1.Range("B1").Copy()
$ws2.Range("A1:C1").PasteSpecial($xlPasteFormats)

$ws2.Columns.Item(1).ColumnWidth = 15.6
$ws2.Columns.Item(2).ColumnWidth = 18.6
$ws2.Columns.Item(3).ColumnWidth = 48.77

$null = $ws2.Range("A2").Select()
$ws2.Activate()
